# Adding Notification Test Case
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Row 18 (TestCase_F17) Result flips from PASS to SKIP now that the newly
# added TestCase_F18 row becomes the last ("PASS") entry.
$ws.Range("E18").Value = "SKIP"

# New TestCase_F18 row.
$ws.Range("A19").Value = "TestCase_F18"
$ws.Range("B19").Value = "OPQA-1099"
$ws.Range("C19").Value = "Verify that Featured Post move down when new notification event occur"
$ws.Range("D19").Value = "Y"
$ws.Range("E19").Value = "PASS"

# Match formatting of the row above (row 18) for the new row.
$ws.Range("A18:E18").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)

# Update selection / view state to match the post-edit workbook.
$ws.Range("A19").Select()
